# Updated cryptos list on Mon Jan 29 09:36:21 UTC 2024 with GitHub Actions
# Refreshes price (col D) and 1h volume-change (col E) figures for the
# crypto table on Sheet1, and swaps the Cosmos/Toncoin row order (rows 29-30).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '42.182.46'
$ws.Range('E2').Value = '  -0.89%  '
$ws.Range('D3').Value = '2.265.04'
$ws.Range('E3').Value = '  -1.13%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.00'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  +0.03%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '307.96'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.24%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '97.44'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +0.50%  '
$ws.Range('E7').Value = '  -1.35%  '
$ws.Range('E8').Value = '  +0.00%  '
$ws.Range('E9').Value = '  -1.51%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '35.24'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -2.96%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0792'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -1.63%  '
$ws.Range('E12').Value = '  +0.39%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '6.80'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +1.10%  '
$ws.Range('D14').Value = '2.615.08'
$ws.Range('E14').Value = '  -1.08%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '14.57'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -0.68%  '
$ws.Range('D16').Value = '2.287.88'
$ws.Range('E16').Value = '  -0.53%  '
$ws.Range('E17').Value = '  -2.05%  '
$ws.Range('D18').Value = '41.987.44'
$ws.Range('E18').Value = '  -1.03%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '12.23'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -4.38%  '
$ws.Range('E20').Value = '  -1.94%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '5.96'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -0.96%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '67.53'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -0.53%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '236.25'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -2.81%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.59'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -0.81%  '
$ws.Range('E25').Value = '  +0.41%  '
$ws.Range('E26').Value = '  -0.08%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '23.56'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -1.93%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '36.58'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -0.03%  '
$ws.Range('B29').Value = 'Cosmos'
$ws.Range('C29').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '9.54'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -0.54%  '
$ws.Range('B30').Value = 'Toncoin'
$ws.Range('C30').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '2.13'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +0.91%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '163.85'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +1.75%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '5.24'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -1.78%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.00'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +0.18%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '3.12'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +0.34%  '
$ws.Range('E35').Value = '  -2.24%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '17.37'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -0.28%  '
$ws.Range('E37').Value = '  +0.19%  '
$ws.Range('E38').Value = '  -3.69%  '
$ws.Range('E39').Value = '  -0.74%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '1.82'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -3.56%  '
$ws.Range('E41').Value = '  -1.19%  '
$ws.Range('E42').Value = '  -6.31%  '
$ws.Range('D43').Value = '1.958.98'
$ws.Range('E43').Value = '  -2.33%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '18.77'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -2.44%  '
$ws.Range('E45').Value = '  -2.01%  '
$ws.Range('E46').Value = '  -3.24%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '9.79'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -4.32%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '53.48'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -1.15%  '
$ws.Range('D49').Value = '2.486.68'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '92.36'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +0.34%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '71.45'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -1.66%  '
